$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update environment data from test21 -> test18
$ws.Range("A2").Value = "https://test18.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test18.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test18.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest18"
$ws.Range("G2").Value = "test18"

# Update the active selection to C8
$ws.Range("C8").Select()
